$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-CellText "D2" "43.149.46"
Set-CellText "E2" "  +2.23%  "
Set-CellText "D3" "2.312.27"
Set-CellText "E3" "  +1.64%  "
Set-CellText "E4" "  +0.04%  "
Set-CellText "D5" "302.38"
Set-CellText "E5" "  +1.41%  "
Set-CellText "D6" "100.31"
Set-CellText "E6" "  +5.43%  "
Set-CellText "E7" "  +2.71%  "
Set-CellText "E8" "  -0.03%  "
Set-CellText "E9" "  +3.59%  "
Set-CellText "E10" "  +3.74%  "
Set-CellText "E11" "  +0.88%  "
Set-CellText "E12" "  +4.01%  "
Set-CellText "D13" "18.09"
Set-CellText "E13" "  +15.00%  "
Set-CellText "D14" "6.86"
Set-CellText "E14" "  +3.40%  "
Set-CellText "D15" "2.671.53"
Set-CellText "E15" "  +1.66%  "
Set-CellText "D16" "2.314.77"
Set-CellText "E16" "  +1.01%  "
Set-CellText "D17" "0.818"
Set-CellText "E17" "  +5.14%  "
Set-CellText "D18" "43.101.73"
Set-CellText "E18" "  +2.18%  "
Set-CellText "D19" "12.61"
Set-CellText "E19" "  +10.13%  "
Set-CellText "D20" "0.0₃0906"
Set-CellText "E20" "  +1.90%  "
Set-CellText "D21" "6.13"
Set-CellText "E22" "  +1.79%  "
Set-CellText "D23" "237.36"
Set-CellText "E23" "  +1.88%  "
Set-CellText "E24" "  +11.50%  "
Set-CellText "E25" "  +0.93%  "
Set-CellText "D26" "0.998"
Set-CellText "E26" "  -0.33%  "
Set-CellText "D27" "24.85"
Set-CellText "E27" "  +4.14%  "
Set-CellText "B28" "Monero"
Set-CellText "C28" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-CellText "D28" "167.91"
Set-CellText "E28" "  +0.68%  "
Set-CellText "B29" "Toncoin"
Set-CellText "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-CellText "D29" "2.09"
Set-CellText "E29" "  -9.16%  "
Set-CellText "D30" "34.11"
Set-CellText "E30" "  +0.83%  "
Set-CellText "D31" "9.18"
Set-CellText "E31" "  +1.56%  "
Set-CellText "E32" "  +0.08%  "
Set-CellText "E33" "  +2.60%  "
Set-CellText "D34" "4.63"
Set-CellText "E34" "  +2.33%  "
Set-CellText "E35" "  +5.01%  "
Set-CellText "D36" "17.09"
Set-CellText "E36" "  +6.19%  "
Set-CellText "D37" "0.0691"
Set-CellText "E37" "  +0.09%  "
Set-CellText "E38" "  +4.06%  "
Set-CellText "E39" "  +4.47%  "
Set-CellText "E40" "  +1.41%  "
Set-CellText "E41" "  +0.76%  "
Set-CellText "E42" "  -1.59%  "
Set-CellText "D43" "2.002.28"
Set-CellText "E43" "  +2.20%  "
Set-CellText "E44" "  +3.74%  "
Set-CellText "E45" "  +5.78%  "
Set-CellText "D46" "17.70"
Set-CellText "E46" "  +1.12%  "
Set-CellText "D47" "2.86"
Set-CellText "E47" "  +2.53%  "
Set-CellText "D48" "55.83"
Set-CellText "E48" "  +7.31%  "
Set-CellText "D49" "2.539.92"
Set-CellText "E49" "  +1.62%  "
Set-CellText "E50" "  +4.91%  "
Set-CellText "E51" "  +1.39%  "
